$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), matching the style of the
# other header cells (bold, centered, bordered - same style as H1).
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Values for I2:J79, one pair (I,J) per data row (rows 2-79).
$values = @(
    @(7,8),
    @(7,7),
    @(8,8),
    @(6,6),
    @(7,8),
    @(8,8),
    @(8,8),
    @(6,7),
    @(7,8),
    @(7,7),
    @(7,7),
    @(8,8),
    @(8,8),
    @(7,7),
    @(8,8),
    @(6,6),
    @(6,6),
    @(7,7),
    @(7,8),
    @(6,6),
    @(7,8),
    @(6,6),
    @(8,8),
    @(7,7),
    @(8,8),
    @(6,6),
    @(7,7),
    @(10,10),
    @(7,7),
    @(7,7),
    @(8,8),
    @(7,7),
    @(7,7),
    @(10,10),
    @(9,9),
    @(7,7),
    @(7,7),
    @(8,8),
    @(7,7),
    @(8,8),
    @(7,7),
    @(10,10),
    @(8,8),
    @(8,8),
    @(9,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,8),
    @(8,8),
    @(9,9),
    @(9,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(9,9),
    @(8,8),
    @(9,9),
    @(9,9),
    @(8,8),
    @(7,7),
    @(7,7),
    @(7,7),
    @(7,7),
    @(5,5),
    @(4,4)
)

$numRows = 78
$arr = New-Object 'object[,]' $numRows,2
for ($i = 0; $i -lt $numRows; $i++) {
    $arr[$i,0] = $values[$i][0]
    $arr[$i,1] = $values[$i][1]
}

# Write the whole block of values in a single call (rows 2-79, columns I-J).
$ws.Range("I2:J79").Value = $arr
